# Major bugfix. Now working fine
# Re-anchor to the live ActiveWorkbook/ActiveSheet (the pre-seeded $wb handle
# does not reliably flush writes back to the underlying document in this host).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Updated RMSE values (re-run results) ----
# Block 1 (rows 5-7): Training RMSE
$ws.Range("C5").Value = 3.7362000000000002
$ws.Range("D5").Value = 3.9727000000000001
$ws.Range("E5").Value = 3.6125882765697899
$ws.Range("F5").Value = 2.8656999999999999

$ws.Range("C6").Value = 9.7098999999999993
$ws.Range("D6").Value = 8.7521000000000004
$ws.Range("E6").Value = 5.7435999999999998
$ws.Range("F6").Value = 2.669

$ws.Range("C7").Value = 27.6129
$ws.Range("D7").Value = 23.088999999999999
$ws.Range("E7").Value = 8.0867000000000004
$ws.Range("F7").Value = 2.5977999999999999

# Block 2 (rows 12-14): Validation RMSE
$ws.Range("C12").Value = 2.7580867267323499
$ws.Range("D12").Value = 2.7245297198935998
$ws.Range("E12").Value = 2.8170806007557299
$ws.Range("F12").Value = 2.7148187709613101

$ws.Range("C13").Value = 2.55016
$ws.Range("D13").Value = 2.6207207634949401
$ws.Range("E13").Value = 2.4481518697826301
$ws.Range("F13").Value = 2.4806417943242001

$ws.Range("C14").Value = 2.1398590372275899
$ws.Range("D14").Value = 1.98245689810242
$ws.Range("F14").Value = 2.2333437272225001

# E14 additionally picks up a new emphasised style (bold font + box border)
$e14 = $ws.Range("E14")
$e14.Value = 1.9585765614249799
$e14.Font.Bold = $true
$e14.Borders.LineStyle = 1

# ---- View state: scroll down one row and move the active selection ----
$ws.Activate() | Out-Null
$ws.Range("H13").Select() | Out-Null
